$d = $word.ActiveDocument

# 1. Insert the "Vi förväntar oss..." paragraph right after the
#    "Nedan presenteras fynd av naturvårdsarter..." paragraph (currently paragraph 3).
$targetText = "Nedan presenteras fynd av naturvårdsarter och fridlysta arter som gjorts i det avverkningsanmälda området, samt relevanta utdrag ur standarderna för FSC, Chain of Custody, Controlled Wood och PEFC. I BILAGA 1 finns artfakta om fridlysta arter."

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq $targetText) {
        $rng = $p.Range
        $rng.Collapse(0)
        $rng.InsertParagraphAfter()
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Text = "Vi förväntar oss att ni återkommer med ett skriftligt svar på vårt klagomål och även beskriver vilka korrigerande åtgärder ni satt in för att rätta till identifierade brister i er efterlevnad av den svenska FSC standarden."
        $found = $true
        break
    }
}

# 2. Remove the old location: two empty paragraphs followed by the
#    "Vi förväntar oss..." paragraph that used to sit right after the
#    "I den avverkningsanmälda skogen ..." comment paragraph (and right
#    before the page-break paragraph).
$commentText = "Kommentar: I den avverkningsanmälda skogen har fridlysta arter sina livsmiljöer och växtplatser. Att skada de fridlysta arternas livsmiljöer, växtplatser eller ekologiska funktion är inte tillåtet enligt artskyddsförordningen"

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq $commentText) {
        $p1 = $d.Paragraphs.Item($i + 1)
        $p2 = $d.Paragraphs.Item($i + 2)
        $p3 = $d.Paragraphs.Item($i + 3)
        if ($p1.Range.Text.Trim().Length -eq 0 -and $p2.Range.Text.Trim().Length -eq 0) {
            $start = $p1.Range.Start
            $end = $p3.Range.End
            $delRange = $d.Range($start, $end)
            $delRange.Delete()
        }
        break
    }
}

# 3. Update the date in the header from 2023-11-13 to 2023-11-14.
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)
    $headers = $sec.Headers
    for ($h = 1; $h -le $headers.Count; $h++) {
        $hdr = $headers.Item($h)
        if ($hdr.Exists) {
            $hdr.Range.Find.Execute("2023-11-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-14", 2) | Out-Null
        }
    }
}
